# Worked on student profile image
# Update the "Student" sheet's Profile Photo (H) and Progress (K) columns:
#   - Profile Photo: "default" -> "default.png" (or a specific uploaded photo
#     filename / full file: URI for the students whose photo was changed)
#   - Progress: numeric fraction (0.4, 0.5, ...) -> text percentage ("40%", "50%", ...)
# Also switches the active sheet/tab back to "Student" with a fresh selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Student")

$profilePhoto = @{
    2  = "default.png"
    3  = "2.png"
    4  = "default.png"
    5  = "default.png"
    6  = "default.png"
    7  = "default.png"
    8  = "default.png"
    9  = "file:/C:/Users/khale/IdeaProjects/NewUniversityManagementSystem/src/main/java/com/example/universitymanagementsystem/ExcelDatabase/studentprofileimages/2.png"
    10 = "default.png"
    11 = "default.png"
}

$progress = @{
    2  = "40%"
    3  = "50%"
    4  = "60%"
    5  = "50%"
    6  = "50%"
    7  = "50%"
    8  = "50%"
    9  = "50%"
    10 = "50%"
    11 = "20%"
}

foreach ($row in 2..11) {
    $ws.Range("H$row").Value2 = $profilePhoto[$row]
    $ws.Range("K$row").Value2 = $progress[$row]
}

# Bring the Student sheet back to the front / make it the active tab, with
# the selection left on H13 (matches the post-edit workbook state).
$ws.Activate()
$ws.Range("H13").Select()
